# feat: add 2022-Q3 data
#
# 1. Insert a brand-new worksheet "2022-Q3" right before "2022-Q2" and
#    populate it with the new quarter's fund holdings.
# 2. Update the "总计" (totals) sheet: insert a new top data row for
#    2022-Q3 and shift every existing quarter's row down by one.
# Every other existing sheet (2022-Q2 .. 2020-Q4) keeps its data exactly
# as-is; only its tab position shifts right by one to make room.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: insert the new "2022-Q3" worksheet before the "2022-Q2" sheet.
# ---------------------------------------------------------------------
$refSheet = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($refSheet)
$q3.Name = "2022-Q3"

# Header row (matches the other quarterly sheets).
$q3.Cells.Item(1,2).Value = "基金代码"
$q3.Cells.Item(1,3).Value = "基金名称"
$q3.Cells.Item(1,4).Value = "基金规模"
$q3.Cells.Item(1,5).Value = "股票总仓位"
$q3.Cells.Item(1,6).Value = "仓位占比"
$q3.Cells.Item(1,7).Value = "持有市值(亿元)"
$q3.Cells.Item(1,8).Value = "仓位排名"

$q3Rows = @(
    @(0, "007216", "浙商中华预期高股息C",               "4.40", "88.55", "7.07", "0.3111", 8),
    @(1, "007178", "浙商中华预期高股息A",               "2.59", "88.55", "7.07", "0.1831", 8),
    @(2, "513690", "博时恒生港股通高股息率ETF",          "3.05", "97.26", "2.49", "0.0759", 7),
    @(3, "159726", "华夏恒生中国内地企业高股息率ETF",     "0.84", "96.48", "2.75", "0.0231", 5),
    @(4, "004532", "民生加银中证港股通高股息精选指数A",   "0.13", "92.87", "5.13", "0.0067", 3),
    @(5, "005702", "恒生前海港股通高股息低波动指数",       "0.20", "94.22", "2.44", "0.0049", 6),
    @(6, "004533", "民生加银中证港股通高股息精选指数C",   "0.08", "92.87", "5.13", "0.0041", 3)
)

foreach ($r in $q3Rows) {
    $rowNum = [int]$r[0] + 2
    # Column A (index) and H (rank) are numeric; B..G are text, mirroring
    # the source workbook's cell types (fund codes keep leading zeros,
    # percentages keep trailing zeros).
    $q3.Cells.Item($rowNum,1).Value = $r[0]
    $q3.Cells.Item($rowNum,2).NumberFormat = "@"
    $q3.Cells.Item($rowNum,2).Value = $r[1]
    $q3.Cells.Item($rowNum,3).NumberFormat = "@"
    $q3.Cells.Item($rowNum,3).Value = $r[2]
    $q3.Cells.Item($rowNum,4).NumberFormat = "@"
    $q3.Cells.Item($rowNum,4).Value = $r[3]
    $q3.Cells.Item($rowNum,5).NumberFormat = "@"
    $q3.Cells.Item($rowNum,5).Value = $r[4]
    $q3.Cells.Item($rowNum,6).NumberFormat = "@"
    $q3.Cells.Item($rowNum,6).Value = $r[5]
    $q3.Cells.Item($rowNum,7).NumberFormat = "@"
    $q3.Cells.Item($rowNum,7).Value = $r[6]
    $q3.Cells.Item($rowNum,8).Value = $r[7]
}

# ---------------------------------------------------------------------
# Step 2: rewrite the "总计" summary sheet with the new 2022-Q3 row on
# top and every other row shifted down by one index/position.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$totalRows = @(
    @(0, "2022-Q3", 7, 0.61),
    @(1, "2022-Q2", 7, 1.2),
    @(2, "2022-Q1", 9, 1.33),
    @(3, "2021-Q4", 3, 0.25),
    @(4, "2021-Q3", 7, 0.26),
    @(5, "2021-Q1", 2, 0.5600000000000001),
    @(6, "2020-Q4", 4, 0.01)
)

foreach ($r in $totalRows) {
    $rowNum = [int]$r[0] + 2
    $total.Cells.Item($rowNum,1).Value = $r[0]
    $total.Cells.Item($rowNum,2).Value = $r[1]
    $total.Cells.Item($rowNum,3).Value = $r[2]
    $total.Cells.Item($rowNum,4).Value = $r[3]
}

# Keep "2020-Q4" (still the last tab) the active/selected sheet, since
# inserting the new worksheet would otherwise steal the tab selection.
$wb.Worksheets.Item("2020-Q4").Activate()

